$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: stash original row 2 (before any changes) into a scratch row
$ws.Range("A2:AY2").Copy($ws.Range("A1000:AY1000"))

# Row 2 <- original row 3
$ws.Range("A2:AY2").ClearContents()
$ws.Range("A3:L3").Copy($ws.Range("A2:L2"))
$ws.Range("N3").Copy($ws.Range("N2"))
$ws.Range("P3:W3").Copy($ws.Range("P2:W2"))
$ws.Range("Y3:AG3").Copy($ws.Range("Y2:AG2"))
$ws.Range("AT3").Copy($ws.Range("AT2"))
$ws.Range("AW3:AY3").Copy($ws.Range("AW2:AY2"))

# Row 3 <- original row 4
$ws.Range("A3:AY3").ClearContents()
$ws.Range("A4:L4").Copy($ws.Range("A3:L3"))
$ws.Range("N4").Copy($ws.Range("N3"))
$ws.Range("P4:W4").Copy($ws.Range("P3:W3"))
$ws.Range("Y4:AG4").Copy($ws.Range("Y3:AG3"))
$ws.Range("AT4").Copy($ws.Range("AT3"))
$ws.Range("AW4:AY4").Copy($ws.Range("AW3:AY3"))

# Row 4 <- original row 5
$ws.Range("A4:AY4").ClearContents()
$ws.Range("A5:L5").Copy($ws.Range("A4:L4"))
$ws.Range("N5").Copy($ws.Range("N4"))
$ws.Range("P5:W5").Copy($ws.Range("P4:W4"))
$ws.Range("Y5:AG5").Copy($ws.Range("Y4:AG4"))
$ws.Range("AT5").Copy($ws.Range("AT4"))
$ws.Range("AW5:AY5").Copy($ws.Range("AW4:AY4"))

# Row 5 <- original row 6
$ws.Range("A5:AY5").ClearContents()
$ws.Range("A6:I6").Copy($ws.Range("A5:I5"))
$ws.Range("P6:W6").Copy($ws.Range("P5:W5"))
$ws.Range("Y6:AB6").Copy($ws.Range("Y5:AB5"))
$ws.Range("AD6:AE6").Copy($ws.Range("AD5:AE5"))
$ws.Range("AG6").Copy($ws.Range("AG5"))
$ws.Range("AT6").Copy($ws.Range("AT5"))
$ws.Range("AW6:AY6").Copy($ws.Range("AW5:AY5"))

# Row 6 <- original row 7
$ws.Range("A6:AY6").ClearContents()
$ws.Range("A7:L7").Copy($ws.Range("A6:L6"))
$ws.Range("N7").Copy($ws.Range("N6"))
$ws.Range("P7:W7").Copy($ws.Range("P6:W6"))
$ws.Range("Y7:AG7").Copy($ws.Range("Y6:AG6"))
$ws.Range("AT7").Copy($ws.Range("AT6"))
$ws.Range("AW7:AY7").Copy($ws.Range("AW6:AY6"))

# Row 7 <- original row 8
$ws.Range("A7:AY7").ClearContents()
$ws.Range("A8:L8").Copy($ws.Range("A7:L7"))
$ws.Range("N8").Copy($ws.Range("N7"))
$ws.Range("P8:W8").Copy($ws.Range("P7:W7"))
$ws.Range("Y8:AG8").Copy($ws.Range("Y7:AG7"))
$ws.Range("AT8").Copy($ws.Range("AT7"))
$ws.Range("AW8:AY8").Copy($ws.Range("AW7:AY7"))

# Row 8 <- original row 9
$ws.Range("A8:AY8").ClearContents()
$ws.Range("A9:I9").Copy($ws.Range("A8:I8"))
$ws.Range("P9:W9").Copy($ws.Range("P8:W8"))
$ws.Range("Y9:AB9").Copy($ws.Range("Y8:AB8"))
$ws.Range("AD9:AE9").Copy($ws.Range("AD8:AE8"))
$ws.Range("AG9").Copy($ws.Range("AG8"))
$ws.Range("AT9").Copy($ws.Range("AT8"))
$ws.Range("AW9:AY9").Copy($ws.Range("AW8:AY8"))

# Row 9 <- original row 10
$ws.Range("A9:AY9").ClearContents()
$ws.Range("A10:I10").Copy($ws.Range("A9:I9"))
$ws.Range("P10:W10").Copy($ws.Range("P9:W9"))
$ws.Range("Y10:AB10").Copy($ws.Range("Y9:AB9"))
$ws.Range("AD10:AE10").Copy($ws.Range("AD9:AE9"))
$ws.Range("AG10").Copy($ws.Range("AG9"))
$ws.Range("AT10").Copy($ws.Range("AT9"))
$ws.Range("AW10:AY10").Copy($ws.Range("AW9:AY9"))

# Row 10 <- original row 11
$ws.Range("A10:AY10").ClearContents()
$ws.Range("A11:L11").Copy($ws.Range("A10:L10"))
$ws.Range("N11").Copy($ws.Range("N10"))
$ws.Range("P11:W11").Copy($ws.Range("P10:W10"))
$ws.Range("Y11:AG11").Copy($ws.Range("Y10:AG10"))
$ws.Range("AT11").Copy($ws.Range("AT10"))
$ws.Range("AW11:AY11").Copy($ws.Range("AW10:AY10"))

# Row 11 <- original row 12
$ws.Range("A11:AY11").ClearContents()
$ws.Range("A12:I12").Copy($ws.Range("A11:I11"))
$ws.Range("P12:W12").Copy($ws.Range("P11:W11"))
$ws.Range("Y12:AB12").Copy($ws.Range("Y11:AB11"))
$ws.Range("AD12:AE12").Copy($ws.Range("AD11:AE11"))
$ws.Range("AG12").Copy($ws.Range("AG11"))
$ws.Range("AT12").Copy($ws.Range("AT11"))
$ws.Range("AW12:AY12").Copy($ws.Range("AW11:AY11"))

# Row 12 <- original row 13
$ws.Range("A12:AY12").ClearContents()
$ws.Range("A13:I13").Copy($ws.Range("A12:I12"))
$ws.Range("P13:W13").Copy($ws.Range("P12:W12"))
$ws.Range("Y13:AB13").Copy($ws.Range("Y12:AB12"))
$ws.Range("AD13:AE13").Copy($ws.Range("AD12:AE12"))
$ws.Range("AG13").Copy($ws.Range("AG12"))
$ws.Range("AT13").Copy($ws.Range("AT12"))
$ws.Range("AW13:AY13").Copy($ws.Range("AW12:AY12"))

# Step 3: move the stashed original row 2 data into row 13
$ws.Range("A13:AY13").ClearContents()
$ws.Range("A1000:I1000").Copy($ws.Range("A13:I13"))
$ws.Range("P1000:W1000").Copy($ws.Range("P13:W13"))
$ws.Range("Y1000:AB1000").Copy($ws.Range("Y13:AB13"))
$ws.Range("AD1000:AE1000").Copy($ws.Range("AD13:AE13"))
$ws.Range("AG1000").Copy($ws.Range("AG13"))
$ws.Range("AT1000").Copy($ws.Range("AT13"))
$ws.Range("AW1000:AY1000").Copy($ws.Range("AW13:AY13"))

# Step 4: clean up the scratch row
$ws.Range("A1000:AY1000").Clear()

Write-Output "done"